# Edit commit: "Mon, Mar 23, 2020  5:06:01 PM"
#
# 1) Swap the deck's colour theme from the custom "Integral" palette to the
#    stock Office 2016+ palette (the file keeps using ppt/theme/theme1.xml,
#    only the colours inside it change).
# 2) Re-point the small two-column table on slide 16 at a different table
#    style (a built-in style GUID rather than the custom "Table_0" style).

$p = $ppt.ActivePresentation

# --- 1) Theme colours: Integral -> Office -------------------------------
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

# Index : Scheme role : Office RGB (as VBA RGB(r,g,b) = r + g*256 + b*65536)
$scheme.Colors(1).RGB  = 0        # dk1      000000
$scheme.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388  # dk2      44546A
$scheme.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407    # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308 # accent5  4472C4
$scheme.Colors(10).RGB = 4697456  # accent6  70AD47
$scheme.Colors(11).RGB = 12673797 # hlink    0563C1
$scheme.Colors(12).RGB = 7491477  # folHlink 954F72

# --- 2) Table style on slide 16 -----------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{8A662AC7-E4A3-49F7-875A-0E644C24A9F9}")
